$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current header values for K1:O1 before moving them:
# K1 = "image name"
# L1 = "library (0 = no, 1 = yes)"
# M1 = "library_base_price"
# N1 = "percentage"
# O1 = "multiplier"
$imageName   = $ws.Range("K1").Value()
$libraryFlag = $ws.Range("L1").Value()
$libraryBase = $ws.Range("M1").Value()
$percentage  = $ws.Range("N1").Value()
$multiplier  = $ws.Range("O1").Value()

# Reorder so percentage and multiplier sit right after prices (column J),
# followed by the former K/L/M (image name, library flag, library_base_price).
$ws.Range("K1").Value = $percentage
$ws.Range("L1").Value = $multiplier
$ws.Range("M1").Value = $imageName
$ws.Range("N1").Value = $libraryFlag
$ws.Range("O1").Value = $libraryBase

# Update the top-left visible cell of the sheet view.
$ws.Application.ActiveWindow.ScrollColumn = 9
